$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the window so row 32 is at the top (matches topLeftCell A31 -> A32)
$excel.Goto($ws.Range("A32"), $true)
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1

# Select the full data range A2:F65 and autofit column A to the SMILES content,
# mirroring the existing best-fit sizing already applied to column D
$ws.Range("A2:F65").Select()
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
